$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Bibi Cell Mundi
$ws.Range("M2").Value = 23495.43
$ws.Range("AG2").Value = 125706.34

# Row 3 - Bibi Cell Vieiralves
$ws.Range("M3").Value = 2350
$ws.Range("AG3").Value = 75669.8

# Row 4 - Bibi Cell Manauara
$ws.Range("L4").Value = 2802
$ws.Range("M4").Value = 2638
$ws.Range("AG4").Value = 39974.89

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("M5").Value = 2702.01
$ws.Range("AG5").Value = 36211.12

# Row 6 - total
$ws.Range("L6").Value = 9776.9
$ws.Range("M6").Value = 31185.44
$ws.Range("AG6").Value = 277562.15
